$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tuần 1 (rows 4-7): "Công việc" / "Kết quả thực hiện" columns
$ws.Range("F4").Value = "Khảo sát hiện trạng"
$ws.Range("G4").Value = "hoàn thành"

$ws.Range("F5").Value = "Vẽ use case "
$ws.Range("G5").Value = "hoàn thành"

$ws.Range("F6").Value = "Tạo link github"
$ws.Range("G6").Value = "hoàn thành"

$ws.Range("F7").Value = "Phát thảo công nghệ"
$ws.Range("G7").Value = "hoàn thành"

# Move the active selection to F17, matching the last cursor position
$ws.Range("F17").Select()
